$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 56949.25
$ws.Range("I6").Value = 1014.5
$ws.Range("J6").Value = 75594.164
$ws.Range("K6").Value = 3043.5
$ws.Range("L6").Value = 226782.492
$ws.Range("M6").Value = -2931.5
$ws.Range("N6").Value = -227006.492
$ws.Range("H43").Value = 1572.0416
$ws.Range("I43").Value = 671.5
$ws.Range("K43").Value = 671.5
$ws.Range("M43").Value = -602.5
$ws.Range("H129").Value = 833.4838999999999
$ws.Range("J129").Value = 1015.46155
$ws.Range("L129").Value = 3046.38465
$ws.Range("N129").Value = -13046.38465
$ws.Range("H138").Value = 2187.1555
$ws.Range("I138").Value = 2246.6875
$ws.Range("J138").Value = 2154.3103
$ws.Range("K138").Value = 6740.0625
$ws.Range("L138").Value = 6462.9309
$ws.Range("M138").Value = -1600.0625
$ws.Range("N138").Value = -16742.9309

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("M3").ClearContents()
$ws.Range("H19").Value = 3000
$ws.Range("I19").Value = 3000
$ws.Range("K19").Value = 3000
$ws.Range("M19").Value = -2771
$ws.Range("H63").Value = 2872.125
$ws.Range("I63").Value = 2568.1428
$ws.Range("K63").Value = 2568.1428
$ws.Range("M63").Value = -1882.1428
$ws.Range("H66").Value = 2872.125
$ws.Range("I66").Value = 2568.1428
$ws.Range("K66").Value = 12840.714
$ws.Range("M66").Value = -9408.714
$ws.Range("H74").Value = 6463665.5
$ws.Range("I74").Value = 12551405
$ws.Range("J74").Value = 55517.844
$ws.Range("K74").Value = 12551405
$ws.Range("L74").Value = 55517.844
$ws.Range("M74").Value = -12550531
$ws.Range("N74").Value = -57265.844
$ws.Range("H77").Value = 6463665.5
$ws.Range("I77").Value = 12551405
$ws.Range("J77").Value = 55517.844
$ws.Range("K77").Value = 62757025
$ws.Range("L77").Value = 277589.22
$ws.Range("M77").Value = -62752657
$ws.Range("N77").Value = -286325.22
$ws.Range("H122").Value = 4832773.5
$ws.Range("I122").Value = 1944.4736
$ws.Range("J122").Value = 27779210
$ws.Range("K122").Value = 5833.4208
$ws.Range("L122").Value = 83337630
$ws.Range("M122").Value = -3383.4208
$ws.Range("N122").Value = -83342530
$ws.Range("H132").Value = 48710.418
$ws.Range("I132").Value = 38796
$ws.Range("J132").Value = 65441
$ws.Range("K132").Value = 116388
$ws.Range("L132").Value = 196323
$ws.Range("M132").Value = -113858
$ws.Range("N132").Value = -201383

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2796.3667
$ws.Range("J134").Value = 3433.0908
$ws.Range("L134").Value = 10299.2724
$ws.Range("N134").Value = -15369.2724

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H12").Value = 992.5
$ws.Range("J12").Value = 0
$ws.Range("L12").Value = 0
$ws.Range("N12").ClearContents()
$ws.Range("H13").Value = 602
$ws.Range("J13").Value = 1000
$ws.Range("L13").Value = 1000
$ws.Range("N13").Value = -1278
$ws.Range("H22").Value = 71429170
$ws.Range("I22").Value = 125000280
$ws.Range("K22").Value = 125000280
$ws.Range("M22").Value = -124999930
$ws.Range("H31").Value = 1843.5454
$ws.Range("I31").Value = 1353.48
$ws.Range("K31").Value = 1353.48
$ws.Range("M31").Value = -1058.48
$ws.Range("H34").Value = 1843.5454
$ws.Range("I34").Value = 1353.48
$ws.Range("K34").Value = 1353.48
$ws.Range("M34").Value = -1151.48
$ws.Range("H132").Value = 24907.883
$ws.Range("I132").Value = 1019.25806
$ws.Range("K132").Value = 3057.77418
$ws.Range("M132").Value = -527.7741799999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 740.1905
$ws.Range("I5").Value = 547.9
$ws.Range("J5").Value = 915
$ws.Range("K5").Value = 1643.7
$ws.Range("L5").Value = 2745
$ws.Range("M5").Value = -1531.7
$ws.Range("N5").Value = -2969
$ws.Range("H76").Value = 3336.3635
$ws.Range("J76").Value = 3570
$ws.Range("L76").Value = 10710
$ws.Range("N76").Value = -11476
$ws.Range("H79").Value = 3336.3635
$ws.Range("J79").Value = 3570
$ws.Range("L79").Value = 10710
$ws.Range("N79").Value = -13362
$ws.Range("H122").Value = 897
$ws.Range("I122").Value = 409.1111
$ws.Range("K122").Value = 3681.9999
$ws.Range("M122").Value = -1231.9999
$ws.Range("H131").Value = 16986.5
$ws.Range("I131").Value = 603.63635
$ws.Range("J131").Value = 20991.2
$ws.Range("K131").Value = 1810.90905
$ws.Range("L131").Value = 62973.60000000001
$ws.Range("M131").Value = 3229.09095
$ws.Range("N131").Value = -73053.60000000001
$ws.Range("H135").Value = 740.1905
$ws.Range("I135").Value = 547.9
$ws.Range("J135").Value = 915
$ws.Range("K135").Value = 4931.099999999999
$ws.Range("L135").Value = 8235
$ws.Range("M135").Value = -2396.099999999999
$ws.Range("N135").Value = -13305

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H7").Value = 40000000
$ws.Range("J7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("N7").ClearContents()
$ws.Range("H8").Value = 40000000
$ws.Range("J8").Value = 0
$ws.Range("L8").Value = 0
$ws.Range("N8").ClearContents()
$ws.Range("H14").Value = 46000052
$ws.Range("I14").Value = 46000052
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 46000052
$ws.Range("L14").Value = 0
$ws.Range("M14").ClearContents()
$ws.Range("N14").Value = -45999884
$ws.Range("H132").Value = 88912.13
$ws.Range("I132").Value = 112327.78
$ws.Range("J132").Value = 73859.21000000001
$ws.Range("K132").Value = 336983.34
$ws.Range("L132").Value = 221577.63
$ws.Range("M132").Value = -334453.34
$ws.Range("N132").Value = -226637.63
$ws.Range("H138").Value = 40000
$ws.Range("J138").Value = 40000
$ws.Range("L138").Value = 40000
$ws.Range("N138").Value = -50280

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 48895.19
$ws.Range("I132").Value = 35213.793
$ws.Range("J132").Value = 79415.234
$ws.Range("K132").Value = 105641.379
$ws.Range("L132").Value = 238245.702
$ws.Range("M132").Value = -103111.379
$ws.Range("N132").Value = -243305.702
$ws.Range("H136").Value = 51400.55
$ws.Range("I136").Value = 29631.914
$ws.Range("K136").Value = 88895.742
$ws.Range("M136").Value = -86345.742
